$wb = $excel.ActiveWorkbook
$wsAbout = $wb.Worksheets.Item("About")
$wsRBFF = $wb.Worksheets.Item("RBFF")

# Data changes: biomass column (G) self-maps instead of mapping to electricity
$wsRBFF.Range("G2").Value = 0
$wsRBFF.Range("G7").Value = 1

# Selection / active sheet changes
$wsAbout.Range("J13").Select() | Out-Null
$wsRBFF.Activate() | Out-Null
$wsRBFF.Range("M7").Select() | Out-Null
